$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.645.02"
$ws.Range("E2").Value = "  +6.47%  "

$ws.Range("D3").Value = "1.772.17"
$ws.Range("E3").Value = "  +3.04%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("E6").Value = "  +3.43%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.95"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.51"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.53%  "

$ws.Range("E10").Value = "  +2.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0659"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").Value = "2.024.91"
$ws.Range("E13").Value = "  +2.98%  "

$ws.Range("D14").Value = "1.773.06"
$ws.Range("E14").Value = "  +3.16%  "

$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").Value = "33.616.84"
$ws.Range("E16").Value = "  +6.41%  "

$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "248.84"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("E21").Value = "  +1.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("E25").Value = "  -2.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.45%  "

$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("E29").Value = "  +1.15%  "

$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0514"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.06%  "

$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("E34").Value = "  +3.22%  "

$ws.Range("E35").Value = "  +3.49%  "

$ws.Range("D36").Value = "1.477.61"

$ws.Range("E37").Value = "  +2.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.626"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "82.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("E41").Value = "  +1.71%  "

$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("E43").Value = "  +3.05%  "

$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0513"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.21%  "

$ws.Range("E46").Value = "  +4.56%  "

$ws.Range("D47").Value = "1.919.60"
$ws.Range("E47").Value = "  +3.31%  "

$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("E49").Value = "  +1.41%  "

$ws.Range("E50").Value = "  +12.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.20%  "
